$d = $word.ActiveDocument
$t = $d.Tables.Item(5)
$newRow = $t.Rows.Add()

$cell1 = $t.Cell(11,1)
$cell1.Range.Text = "19/09/16"

$cell2 = $t.Cell(11,2)
Write-Output "cell2 range start=$($cell2.Range.Start) end=$($cell2.Range.End) text=[$($cell2.Range.Text)]"

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr><w:spacing w:after="0"/></w:pPr>
  <w:r><w:t>Publicação de novas versões d</w:t></w:r>
  <w:r><w:t>os casos de usos:</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:spacing w:after="0"/></w:pPr>
  <w:r><w:t xml:space="preserve">ARRUC0210 - Gerar </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>DARE-e</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> e do </w:t></w:r>
  <w:r><w:t>modelo de casos de usos</w:t></w:r>
  <w:r><w:t xml:space="preserve"> para ajuste na regra de negócio ARRRN0201</w:t></w:r>
  <w:r><w:t>.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:spacing w:after="0"/></w:pPr>
  <w:r><w:t xml:space="preserve">Considerado acréscimo de mais </w:t></w:r>
  <w:r><w:t>25</w:t></w:r>
  <w:r><w:t xml:space="preserve"> dias no prazo de execução para OS por causa do tempo que levou pa</w:t></w:r>
  <w:r><w:t>ra responder os Mantis 138</w:t></w:r>
  <w:r><w:t>.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$cell2.Range.InsertXML($xml)

$cell2c = $t.Cell(11,2)
Write-Output "cell2 final text=[$($cell2c.Range.Text)]"
